$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 118.77778
$ws.Cells.Item(9, 9).Value = 121.125
$ws.Cells.Item(9, 11).Value = 121.125
$ws.Cells.Item(9, 13).Value = 47.875
$ws.Cells.Item(19, 8).Value = 8413.959999999999
$ws.Cells.Item(19, 10).Value = 13762.533
$ws.Cells.Item(19, 12).Value = 13762.533
$ws.Cells.Item(19, 14).Value = -14112.533
$ws.Cells.Item(98, 8).Value = 1283.2142
$ws.Cells.Item(98, 9).Value = 1283.2142
$ws.Cells.Item(98, 11).Value = 1283.2142
$ws.Cells.Item(98, 13).Value = 214.7858000000001
$ws.Cells.Item(122, 8).Value = 1283.2142
$ws.Cells.Item(122, 9).Value = 1283.2142
$ws.Cells.Item(122, 11).Value = 3849.6426
$ws.Cells.Item(122, 13).Value = -1399.6426
$ws.Cells.Item(132, 8).Value = 2215.48
$ws.Cells.Item(132, 9).Value = 2016.4894
$ws.Cells.Item(132, 11).Value = 6049.468199999999
$ws.Cells.Item(132, 13).Value = -3519.468199999999
$ws.Cells.Item(135, 8).Value = 12821569
$ws.Cells.Item(135, 9).Value = 19231650
$ws.Cells.Item(135, 10).Value = 1408.2307
$ws.Cells.Item(135, 11).Value = 173084850
$ws.Cells.Item(135, 12).Value = 12674.0763
$ws.Cells.Item(135, 13).Value = -173082315
$ws.Cells.Item(135, 14).Value = -17744.0763
$ws.Cells.Item(137, 8).Value = 1489496.9
$ws.Cells.Item(137, 9).Value = 2382095.8
$ws.Cells.Item(137, 10).Value = 1831.9048
$ws.Cells.Item(137, 11).Value = 7146287.399999999
$ws.Cells.Item(137, 12).Value = 5495.7144
$ws.Cells.Item(137, 13).Value = -7143737.399999999
$ws.Cells.Item(137, 14).Value = -10595.7144
$ws.Cells.Item(138, 8).Value = 4149.1216
$ws.Cells.Item(138, 9).Value = 4401.933
$ws.Cells.Item(138, 10).Value = 4084.8474
$ws.Cells.Item(138, 11).Value = 13205.799
$ws.Cells.Item(138, 12).Value = 12254.5422
$ws.Cells.Item(138, 13).Value = -8065.798999999999
$ws.Cells.Item(138, 14).Value = -22534.5422

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 13833.1875
$ws.Cells.Item(31, 10).Value = 44004.75
$ws.Cells.Item(31, 12).Value = 44004.75
$ws.Cells.Item(31, 14).Value = -44592.75
$ws.Cells.Item(74, 8).Value = 12501726
$ws.Cells.Item(74, 9).Value = 856
$ws.Cells.Item(74, 10).Value = 33336510
$ws.Cells.Item(74, 11).Value = 856
$ws.Cells.Item(74, 12).Value = 33336510
$ws.Cells.Item(74, 13).Value = 18
$ws.Cells.Item(74, 14).Value = -33338258
$ws.Cells.Item(77, 8).Value = 12501726
$ws.Cells.Item(77, 9).Value = 856
$ws.Cells.Item(77, 10).Value = 33336510
$ws.Cells.Item(77, 11).Value = 4280
$ws.Cells.Item(77, 12).Value = 166682550
$ws.Cells.Item(77, 13).Value = 88
$ws.Cells.Item(77, 14).Value = -166691286
$ws.Cells.Item(93, 8).Value = 64482.668
$ws.Cells.Item(93, 10).Value = 64482.668
$ws.Cells.Item(93, 12).Value = 64482.668
$ws.Cells.Item(93, 14).Value = -69474.66800000001
$ws.Cells.Item(132, 8).Value = 1482628.4
$ws.Cells.Item(132, 9).Value = 2531.257
$ws.Cells.Item(132, 10).Value = 4529887
$ws.Cells.Item(132, 11).Value = 7593.771000000001
$ws.Cells.Item(132, 12).Value = 13589661
$ws.Cells.Item(132, 13).Value = -5063.771000000001
$ws.Cells.Item(132, 14).Value = -13594721

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 408.46155
$ws.Cells.Item(94, 9).Value = 408.46155
$ws.Cells.Item(94, 11).Value = 408.46155
$ws.Cells.Item(94, 13).Value = 42.53845000000001
$ws.Cells.Item(134, 8).Value = 2934.5144
$ws.Cells.Item(134, 9).Value = 2955.2
$ws.Cells.Item(134, 11).Value = 8865.599999999999
$ws.Cells.Item(134, 13).Value = -6330.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 58.444443
$ws.Cells.Item(7, 10).Value = 93.333336
$ws.Cells.Item(7, 12).Value = 93.333336
$ws.Cells.Item(7, 14).Value = -319.333336
$ws.Cells.Item(31, 8).Value = 5270.2856
$ws.Cells.Item(31, 9).Value = 1483.0605
$ws.Cells.Item(31, 10).Value = 7425.0864
$ws.Cells.Item(31, 11).Value = 1483.0605
$ws.Cells.Item(31, 12).Value = 7425.0864
$ws.Cells.Item(31, 13).Value = -1188.0605
$ws.Cells.Item(31, 14).Value = -8015.0864
$ws.Cells.Item(34, 8).Value = 5270.2856
$ws.Cells.Item(34, 9).Value = 1483.0605
$ws.Cells.Item(34, 10).Value = 7425.0864
$ws.Cells.Item(34, 11).Value = 1483.0605
$ws.Cells.Item(34, 12).Value = 7425.0864
$ws.Cells.Item(34, 13).Value = -1281.0605
$ws.Cells.Item(34, 14).Value = -7829.0864
$ws.Cells.Item(58, 8).Value = 1463.0834
$ws.Cells.Item(58, 9).Value = 1079
$ws.Cells.Item(58, 11).Value = 1079
$ws.Cells.Item(58, 13).Value = -876
$ws.Cells.Item(132, 8).Value = 12580314
$ws.Cells.Item(132, 9).Value = 12821889
$ws.Cells.Item(132, 11).Value = 38465667
$ws.Cells.Item(132, 13).Value = -38463137
$ws.Cells.Item(136, 8).Value = 1463.0834
$ws.Cells.Item(136, 9).Value = 1079
$ws.Cells.Item(136, 11).Value = 3237
$ws.Cells.Item(136, 13).Value = -687

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 577.7727
$ws.Cells.Item(7, 10).Value = 689.2857
$ws.Cells.Item(7, 12).Value = 2067.8571
$ws.Cells.Item(7, 14).Value = -2291.8571
$ws.Cells.Item(80, 8).Value = 4499.421
$ws.Cells.Item(80, 9).Value = 3734
$ws.Cells.Item(80, 10).Value = 4642.9375
$ws.Cells.Item(80, 11).Value = 11202
$ws.Cells.Item(80, 12).Value = 13928.8125
$ws.Cells.Item(80, 13).Value = -10266
$ws.Cells.Item(80, 14).Value = -15800.8125
$ws.Cells.Item(83, 8).Value = 4499.421
$ws.Cells.Item(83, 9).Value = 3734
$ws.Cells.Item(83, 10).Value = 4642.9375
$ws.Cells.Item(83, 11).Value = 33606
$ws.Cells.Item(83, 12).Value = 41786.4375
$ws.Cells.Item(83, 13).Value = -28926
$ws.Cells.Item(83, 14).Value = -51146.4375
$ws.Cells.Item(113, 8).Value = 511.36508
$ws.Cells.Item(113, 10).Value = 568.1
$ws.Cells.Item(113, 12).Value = 1704.3
$ws.Cells.Item(113, 14).Value = -6044.3
$ws.Cells.Item(129, 8).Value = 674588.5600000001
$ws.Cells.Item(129, 9).Value = 448.6154
$ws.Cells.Item(129, 10).Value = 948457.9
$ws.Cells.Item(129, 11).Value = 1345.8462
$ws.Cells.Item(129, 12).Value = 2845373.7
$ws.Cells.Item(129, 13).Value = 3654.1538
$ws.Cells.Item(129, 14).Value = -2855373.7
$ws.Cells.Item(132, 8).Value = 2678.9487
$ws.Cells.Item(132, 9).Value = 2549.087
$ws.Cells.Item(132, 11).Value = 22941.783
$ws.Cells.Item(132, 13).Value = -20411.783

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(68, 8).Value = 98295
$ws.Cells.Item(68, 10).Value = 98295
$ws.Cells.Item(68, 12).Value = 98295
$ws.Cells.Item(68, 14).Value = -99917
$ws.Cells.Item(71, 8).Value = 98295
$ws.Cells.Item(71, 10).Value = 98295
$ws.Cells.Item(71, 12).Value = 294885
$ws.Cells.Item(71, 14).Value = -302997
$ws.Cells.Item(80, 8).Value = 19565064
$ws.Cells.Item(80, 9).Value = 58556556
$ws.Cells.Item(80, 10).Value = 69316.664
$ws.Cells.Item(80, 11).Value = 58556556
$ws.Cells.Item(80, 12).Value = 69316.664
$ws.Cells.Item(80, 13).Value = -58555558
$ws.Cells.Item(80, 14).Value = -71312.664
$ws.Cells.Item(83, 8).Value = 19565064
$ws.Cells.Item(83, 9).Value = 58556556
$ws.Cells.Item(83, 10).Value = 69316.664
$ws.Cells.Item(83, 11).Value = 292782780
$ws.Cells.Item(83, 12).Value = 346583.32
$ws.Cells.Item(83, 13).Value = -292777788
$ws.Cells.Item(83, 14).Value = -356567.32
$ws.Cells.Item(113, 8).Value = 62875
$ws.Cells.Item(113, 9).Value = 80166
$ws.Cells.Item(113, 10).Value = 2356.5
$ws.Cells.Item(113, 11).Value = 80166
$ws.Cells.Item(113, 12).Value = 2356.5
$ws.Cells.Item(113, 13).Value = -77996
$ws.Cells.Item(113, 14).Value = -6696.5
$ws.Cells.Item(126, 8).Value = 10000
$ws.Cells.Item(126, 9).Value = 10000
$ws.Cells.Item(126, 11).Value = 30000
$ws.Cells.Item(126, 13).Value = -27530
$ws.Cells.Item(132, 8).Value = 47626930
$ws.Cells.Item(132, 9).Value = 71438104
$ws.Cells.Item(132, 10).Value = 4575
$ws.Cells.Item(132, 11).Value = 214314312
$ws.Cells.Item(132, 12).Value = 13725
$ws.Cells.Item(132, 13).Value = -214311782
$ws.Cells.Item(132, 14).Value = -18785

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6423.6665
$ws.Cells.Item(22, 9).Value = 735.7143
$ws.Cells.Item(22, 10).Value = 10043.272
$ws.Cells.Item(22, 11).Value = 735.7143
$ws.Cells.Item(22, 12).Value = 10043.272
$ws.Cells.Item(22, 13).Value = -440.7143
$ws.Cells.Item(22, 14).Value = -10633.272
$ws.Cells.Item(27, 8).Value = 6423.6665
$ws.Cells.Item(27, 9).Value = 735.7143
$ws.Cells.Item(27, 10).Value = 10043.272
$ws.Cells.Item(27, 11).Value = 735.7143
$ws.Cells.Item(27, 12).Value = 10043.272
$ws.Cells.Item(27, 13).Value = -628.7143
$ws.Cells.Item(27, 14).Value = -10257.272
$ws.Cells.Item(61, 8).Value = 5211.5557
$ws.Cells.Item(61, 9).Value = 2849.75
$ws.Cells.Item(61, 11).Value = 2849.75
$ws.Cells.Item(61, 13).Value = -2647.75
$ws.Cells.Item(100, 8).Value = 78897.836
$ws.Cells.Item(100, 9).Value = 153326.67
$ws.Cells.Item(100, 10).Value = 4469
$ws.Cells.Item(100, 11).Value = 153326.67
$ws.Cells.Item(100, 12).Value = 4469
$ws.Cells.Item(100, 13).Value = -152785.67
$ws.Cells.Item(100, 14).Value = -5551
$ws.Cells.Item(113, 8).Value = 5211.5557
$ws.Cells.Item(113, 9).Value = 2849.75
$ws.Cells.Item(113, 11).Value = 2849.75
$ws.Cells.Item(113, 13).Value = -679.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 6749
$ws.Cells.Item(7, 10).Value = 6749
$ws.Cells.Item(7, 12).Value = 6749
$ws.Cells.Item(7, 14).Value = -6975
$ws.Cells.Item(100, 8).Value = 536.4
$ws.Cells.Item(100, 9).Value = 526.6923
$ws.Cells.Item(100, 10).Value = 554.4286
$ws.Cells.Item(100, 11).Value = 1053.3846
$ws.Cells.Item(100, 12).Value = 1108.8572
$ws.Cells.Item(100, 13).Value = -512.3846000000001
$ws.Cells.Item(100, 14).Value = -2190.8572
$ws.Cells.Item(132, 8).Value = 3789834
$ws.Cells.Item(132, 9).Value = 1807.0878
$ws.Cells.Item(132, 11).Value = 5421.2634
$ws.Cells.Item(132, 13).Value = -2891.2634
$ws.Cells.Item(137, 8).Value = 49916.668
$ws.Cells.Item(137, 10).Value = 49916.668
$ws.Cells.Item(137, 12).Value = 49916.668
$ws.Cells.Item(137, 14).Value = -60116.668
